$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status moved from "In Translation" to "Ready for handoff" everywhere it appears.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Timestamps bumped forward by the new handoff generation run.
$overview.Range("G2").Value = "2016-08-12 21:12:11"
$dede.Range("H2").Value = "2016-08-12 21:12:11"
$zhcn.Range("H2").Value = "2016-08-12 21:12:00"

# The longer "Ready for handoff" text widens the Status columns (re-fit).
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
